# Append a new data row (row 70) to the bottom of the daily log sheet,
# matching the "2025/10/06 月 21 64" entry added upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 70

# Column A holds a date formatted as plain text (e.g. "2025/09/22") in every
# existing row, not a real Excel date. Force text formatting on the target
# cell first so the "2025/10/06" string isn't auto-coerced into a date
# serial number when assigned.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025/10/06"

$ws.Cells.Item($row, 2).Value = "月"
$ws.Cells.Item($row, 3).Value = 21
$ws.Cells.Item($row, 4).Value = 64
